$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# --- Section 1: "   * 12121q..." through "   * warwgg" bullet list text tweaks ---
Replace-InParagraph 15 "12121qfasafsvsas vavawvasvsav" "12121qfasafsvsas vavawvasvsavssst,tty"
Replace-InParagraph 16 "MW" "MWyy"
Replace-InParagraph 17 "On Quotation (M&P Below Lacks)" "RSP ITEM"
Replace-InParagraph 18 "sdga" "sdgaty,"
Replace-InParagraph 19 "warwgg" "warwggty,ty,"

# --- Section 2: "None"/"None"/"17-07-2025"/"17-07-2025" -> new dates ---
Replace-InParagraph 21 "None" "09-07-2025"
Replace-InParagraph 22 "None" "QERGG"
Replace-InParagraph 23 "17-07-2025" "02-07-2025"
Replace-InParagraph 24 "17-07-2025" "09-07-2025"

# --- Section 3: add <w:lastRenderedPageBreak/> before "Special Placeholders" run ---
Replace-ParagraphXml 42 '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">  Special Placeholders</w:t></w:r>'

# --- Section 4: Date placeholder examples ---
Replace-InParagraph 53 "26-07-2025" "28-07-2025"
Replace-InParagraph 54 "07-26-2025" "07-28-2025"

# --- Section 5: replace "2025-07-26" paragraph with split runs + proofErr markers ---
Replace-ParagraphXml 55 '<w:r><w:t xml:space="preserve">   * [</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DATE:YYYY</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-MM-DD]</w:t></w:r>'

# --- Section 6: insert three empty paragraphs at the very beginning of the body ---
$beginXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$beginRange = $d.Range(0, 0)
$beginRange.InsertXML($beginXml)
